# Reorder slides: move the slide currently in position 5
# (the "AI EXAMPLE" slide) so that it becomes slide 3, shifting
# the slides that were in positions 3 and 4 back by one.
$p = $ppt.ActivePresentation
$p.Slides.Item(5).MoveTo(3)
